$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the updated "Price" cells are digit-only strings (e.g. "0.999")
# that Excel would otherwise auto-convert to numbers on assignment. The source
# data keeps these as plain text, so mark them as Text before writing the new
# values. Areas are formatted individually (looping $range.Areas) because
# setting .NumberFormat directly on a multi-area Union only reliably applies
# to its first area; looping keeps every cell correct while still reusing a
# single shared style.
$textCells = $excel.Union($ws.Range("D4"), $ws.Range("D5"), $ws.Range("D6"), $ws.Range("D9"), $ws.Range("D10"), $ws.Range("D11"), $ws.Range("D12"), $ws.Range("D13"), $ws.Range("D16"), $ws.Range("D17"), $ws.Range("D19"), $ws.Range("D20"), $ws.Range("D22"), $ws.Range("D23"), $ws.Range("D26"), $ws.Range("D28"), $ws.Range("D29"), $ws.Range("D31"), $ws.Range("D32"), $ws.Range("D33"), $ws.Range("D34"), $ws.Range("D35"), $ws.Range("D36"), $ws.Range("D37"), $ws.Range("D41"), $ws.Range("D45"), $ws.Range("D47"), $ws.Range("D48"), $ws.Range("D49"), $ws.Range("D50"))
foreach ($area in $textCells.Areas) {
    $area.NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "42.943.77"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "2.538.68"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "316.97"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "97.87"
$ws.Range("E6").Value = "  +2.73%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "0.536"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").Value = "36.06"
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("D11").Value = "0.0813"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "7.63"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").Value = "0.111"
$ws.Range("E13").Value = "  -2.48%  "
$ws.Range("D14").Value = "2.922.30"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "2.547.25"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").Value = "15.17"
$ws.Range("E16").Value = "  -2.84%  "
$ws.Range("D17").Value = "0.851"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "42.957.09"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").Value = "6.89"
$ws.Range("E19").Value = "  +5.28%  "
$ws.Range("D20").Value = "12.78"
$ws.Range("E20").Value = "  -2.71%  "
$ws.Range("D21").Value = "0.0₃0965"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").Value = "69.70"
$ws.Range("E22").Value = "  -2.23%  "
$ws.Range("D23").Value = "253.23"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "26.45"
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").Value = "2.40"
$ws.Range("E28").Value = "  +2.72%  "
$ws.Range("D29").Value = "41.10"
$ws.Range("E29").Value = "  +4.68%  "
$ws.Range("E30").Value = "  +3.64%  "
$ws.Range("D31").Value = "5.93"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").Value = "157.52"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").Value = "2.17"
$ws.Range("E33").Value = "  +3.55%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "3.36"
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "2.71"
$ws.Range("E35").Value = "  +4.10%  "
$ws.Range("D36").Value = "19.07"
$ws.Range("E36").Value = "  -3.73%  "
$ws.Range("D37").Value = "0.0789"
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("E39").Value = "  +17.62%  "
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").Value = "21.91"
$ws.Range("E41").Value = "  -9.78%  "
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("D45").Value = "3.31"
$ws.Range("E45").Value = "  -2.19%  "
$ws.Range("D46").Value = "2.017.81"
$ws.Range("E46").Value = "  -2.37%  "
$ws.Range("D47").Value = "9.10"
$ws.Range("E47").Value = "  +3.14%  "
$ws.Range("D48").Value = "84.49"
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("D49").Value = "76.73"
$ws.Range("E49").Value = "  +2.78%  "
$ws.Range("D50").Value = "106.35"
$ws.Range("E50").Value = "  +4.69%  "
$ws.Range("D51").Value = "2.775.15"
$ws.Range("E51").Value = "  +0.46%  "
